$wb = $excel.ActiveWorkbook

# --- Sheet "org" (sheet1): new row 7 ---
$wsOrg = $wb.Worksheets.Item("org")
$wsOrg.Range("A7").Value = "test_03"
$wsOrg.Range("B7").Value = "deletingOrganization"
$wsOrg.Range("C7").Value = "SkillRary"
$wsOrg.Range("D7").Value = "No Organization Found !"
$wsOrg.Range("D7").Font.Name = "Arial"
$wsOrg.Range("D7").Font.Color = 0

# --- Sheet "contact" (sheet2): new row 3 ---
$wsContact = $wb.Worksheets.Item("contact")
$wsContact.Range("B3").Value = "deletingContWithOrg"
$wsContact.Range("C3").Value = "SkillRary"
$wsContact.Range("D3").Value = "Customer"
$wsContact.Range("E3").Value = "Banking"
$wsContact.Range("F3").Value = "deepak"
$wsContact.Range("G3").Value = "No Contact Found !"
$wsContact.Range("G3").Font.Name = "Arial"
$wsContact.Range("G3").Font.Color = 0
